# Furnish Developer Guide with Design Considerations.
#
# The "Straight Arrow Connector 57" (id 155) that links the two
# "Rectangle 11" shapes (id 152 "EventsCenter" and id 153
# "ShowHelp / RequestEvent") is re-glued: its start end is attached to
# the right-hand connection site of shape 152 and its tail end to the
# left-hand connection site of shape 153. Re-attaching the connector
# also repositions/resizes it (PowerPoint recalculates the straight
# line between the two connection sites) and drops the now-unneeded
# horizontal flip while keeping the vertical one.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connector  = $s.Shapes.Item(22)  # Straight Arrow Connector 57 (id 155)
$startShape = $s.Shapes.Item(20)  # Rectangle 11 "EventsCenter" (id 152)
$endShape   = $s.Shapes.Item(21)  # Rectangle 11 "ShowHelp/RequestEvent" (id 153)

# Re-glue the connector ends to the two rectangles (connection site
# index 1 = right side of shape 152, index 3 = left side of shape 153).
$connector.ConnectorFormat.BeginConnect($startShape, 1)
$connector.ConnectorFormat.EndConnect($endShape, 3)

# Drive the resulting geometry explicitly (mirrors what PowerPoint's
# glue-point recalculation produces once the connector is re-attached):
# the connector no longer needs to be flipped horizontally, still needs
# to be flipped vertically, and its bounding box shrinks from the old
# off=(1529913,3186327) ext=(682330,3489) to the new
# off=(1537309,3169661) ext=(622047,7766), in EMU.
$connector.HorizontalFlip = 0
$connector.VerticalFlip = 1

$connector.Left   = 121.04795465590551
$connector.Top    = 249.5796127992126
$connector.Width  = 48.98007874015748
$connector.Height = 0.611496062992126
